$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ATR-I.1.3")

# The "agrario" sector label (rows 27-45, column A) was renamed to "agricultura".
$rng = $ws.Range("A27:A45")
$rng.Value = "agricultura"

# Reflect the selection/view state shown in the saved file: the user had
# selected A27:A45 (the edited range) with the frozen pane scrolled to row 38.
$ws.Range("A27:A45").Select()
